# Update the LR-pair (Efna5-Epha3) data table with new TPM-derived values.
#
# - The shared string formerly read by the worksheet as "Resolving-Mac" is
#   renamed to "Inflammatory-Mac" (every former reference to it is updated).
# - Several rows (4/5, 8/9, 12/13) swap their "Target cluster" (column D)
#   between "MuSCs" and the renamed "Inflammatory-Mac" label.
# - Columns E:T (per-row/per-pair TPM statistics) are refreshed with the
#   values recomputed from the new TPM input.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @{ D = "ECs"; E = 1; F = 0.3333333333333333; G = 0.2457683333333333; H = 0.737305; I = 0.1447271191911903; J = 0.1575855905380038; K = 2; L = 0.6666666666666666; M = 0.107177; N = 0.321531; O = 0.003526763356587491; P = 0.003549676734010809; Q = 0.02634071266166667; R = 0.237066413955; S = 0.0005104183006679603; T = 0.0005593779043481057 }
    3 = @{ D = "FAPs"; E = 1; F = 0.3333333333333333; G = 0.2457683333333333; H = 0.737305; I = 0.1447271191911903; J = 0.1575855905380038; K = 3; L = 1; M = 29.65321533333333; N = 88.95964599999999; O = 0.9757678722356318; P = 0.9821074349659524; Q = 7.287821310447778; R = 65.59039179403; S = 0.1412200731479805; T = 0.1547659801108737 }
    4 = @{ D = "Inflammatory-Mac"; E = 1; F = 0.3333333333333333; G = 0.2457683333333333; H = 0.737305; I = 0.1447271191911903; J = 0.1575855905380038; K = 2; L = 0.6666666666666666; M = 0.04072766666666667; N = 0.122183; O = 0.001340183457265176; P = 0.001348890627627329; Q = 0.01000957075722222; R = 0.090086136815; S = 0.0001939608909576787; T = 0.0002125657261258311 }
    5 = @{ D = "MuSCs"; E = 1; F = 0.3333333333333333; G = 0.2457683333333333; H = 0.737305; I = 0.1447271191911903; J = 0.1575855905380038; K = 2; L = 1; M = 0.5885005; N = 1.177001; O = 0.01936518095051565; P = 0.01299399767240936; Q = 0.1446347870508333; R = 0.8678087223049999; S = 0.002802666851584248; T = 0.002047666796656076 }
    6 = @{ D = "ECs"; E = 3; F = 1; G = 1.036691; H = 3.110073; I = 0.6104826439049008; J = 0.6647217777192627; K = 2; L = 0.6666666666666666; M = 0.107177; N = 0.321531; O = 0.003526763356587491; P = 0.003549676734010809; Q = 0.111109431307; R = 0.999984881763; S = 0.002153027818356454; T = 0.002359547428960371 }
    7 = @{ D = "FAPs"; E = 3; F = 1; G = 1.036691; H = 3.110073; I = 0.6104826439049008; J = 0.6647217777192627; K = 3; L = 1; M = 29.65321533333333; N = 88.95964599999999; O = 0.9757678722356318; P = 0.9821074349659524; Q = 30.74122145712867; R = 276.670993114158; S = 0.5956893504798679; T = 0.652828200081873 }
    8 = @{ D = "Inflammatory-Mac"; E = 3; F = 1; G = 1.036691; H = 3.110073; I = 0.6104826439049008; J = 0.6647217777192627; K = 2; L = 0.6666666666666666; M = 0.04072766666666667; N = 0.122183; O = 0.001340183457265176; P = 0.001348890627627329; Q = 0.04222200548433334; R = 0.379998049359; S = 0.0008181587403088555; T = 0.0008966369759452899 }
    9 = @{ D = "MuSCs"; E = 3; F = 1; G = 1.036691; H = 3.110073; I = 0.6104826439049008; J = 0.6647217777192627; K = 2; L = 1; M = 0.5885005; N = 1.177001; O = 0.01936518095051565; P = 0.01299399767240936; Q = 0.6100931718455; R = 3.660559031073; S = 0.01182210686636762; T = 0.008637393232483914 }
    10 = @{ D = "ECs"; E = 2; F = 1; G = 0.4156905; H = 0.831381; I = 0.2447902369039089; J = 0.1776926317427335; K = 2; L = 0.6666666666666666; M = 0.107177; N = 0.321531; O = 0.003526763356587491; P = 0.003549676734010809; Q = 0.0445524607185; R = 0.267314764311; S = 0.0008633172375630767; T = 0.0006307514007023315 }
    11 = @{ D = "FAPs"; E = 2; F = 1; G = 0.4156905; H = 0.831381; I = 0.2447902369039089; J = 0.1776926317427335; K = 3; L = 1; M = 29.65321533333333; N = 88.95964599999999; O = 0.9757678722356318; P = 0.9821074349659524; Q = 12.326559908521; R = 73.95935945112599; S = 0.2388584486077834; T = 0.1745132547732056 }
    12 = @{ D = "Inflammatory-Mac"; E = 2; F = 1; G = 0.4156905; H = 0.831381; I = 0.2447902369039089; J = 0.1776926317427335; K = 2; L = 0.6666666666666666; M = 0.04072766666666667; N = 0.122183; O = 0.001340183457265176; P = 0.001348890627627329; Q = 0.0169301041205; R = 0.101580624723; S = 0.0003280638259986421; T = 0.0002396879255562076 }
    13 = @{ D = "MuSCs"; E = 2; F = 1; G = 0.4156905; H = 0.831381; I = 0.2447902369039089; J = 0.1776926317427335; K = 2; L = 1; M = 0.5885005; N = 1.177001; O = 0.01936518095051565; P = 0.01299399767240936; Q = 0.24463406709525; R = 0.9785362683810001; S = 0.00474040723256379; T = 0.002308937643269373 }
}

foreach ($r in $rowData.Keys) {
    $values = $rowData[$r]
    foreach ($col in $values.Keys) {
        $ref = "$col$r"
        $ws.Range($ref).Value = $values[$col]
    }
}
